$d = $word.ActiveDocument

# --- Paragraph 6: Objetivos (PT) -> becomes "Programa resumido" PT text ---
$p = $d.Paragraphs.Item(6)
$full = $d.Range($p.Range.Start, $p.Range.End - 1)
$full.Text = "Propriedade gerais dos compostos orgânicos. Estrutura, métodos de obtenção, propriedades físicas, reações dos hidrocarbonetos alifáticos e aromáticos, haletos orgânicos, álcoois e características estruturais como Estereoquímica e a relação estrutura-reatividade."

# --- Paragraph 7: Objetivos (EN, italic) -> becomes "Programa resumido" EN text ---
$p = $d.Paragraphs.Item(7)
$full = $d.Range($p.Range.Start, $p.Range.End - 1)
$full.Text = "General property of organic compounds. Physical properties, reactions of aliphatic and aromatic hydrocarbons, organic halides, ethers, alcohols and structural characteristics as sstereochemistry and structure-reactivity."
$full.Font.Italic = 1

# --- Paragraph 9: Docente ListBullet ("210064...") -> becomes Objetivos PT text ---
$p = $d.Paragraphs.Item(9)
$full = $d.Range($p.Range.Start, $p.Range.End - 1)
$full.Text = "Gerais - Apresentar e Ensinar conceitos relacionados com o desenvolvimento de dispositivos tecnológicos úteis para a aplicação como sensores, geradores de energia e catálise. Abordar problemáticas sociais e ambientais com as quais a engenharia química e o desenvolvimento de novas tecnologias estão relacionados.`v`vEspecíficos – Compreender e descrever o mecanismo das reações orgânicas e a sua importância para o aprimoramento e desenvolvimento de processos industriais sintéticos e de etapas de formulação. Aprofundar o conceito de estrutura-reatividade e propriedades dos materiais."

# --- Paragraph 11: Programa resumido PT -> becomes Programa PT text ---
$p = $d.Paragraphs.Item(11)
$full = $d.Range($p.Range.Start, $p.Range.End - 1)
$full.Text = "1. Estrutura e propriedades fundamentadas em grupos funcionais. `v2. Compostos orgânicos utilizados em materiais: classificação e aplicação.`v3. Conceitos físico-químicos relacionados às propriedades.`v4. Moléculas orgânicas na formação de Cristais, Géis, Associações Supramoleculares.`v3. Processo do estado fundamental, excitado e eventos de oxido-redução.`v5. Técnicas avançadas de caracterização."

# --- Paragraph 12: Programa resumido EN (italic) -> becomes Objetivos EN text ---
$p = $d.Paragraphs.Item(12)
$full = $d.Range($p.Range.Start, $p.Range.End - 1)
$full.Text = "Overview - Introduce and teach concepts of organic chemistry as important tools for understanding strategies and industrial and technological operations. Address social and environmental issues with which chemical engineering is related, making them thus able to exercise Chemical Engineer function, and realize the changes that are necessary.`v`v`v`vSpecific - Understand and describe the mechanism of organic reactions and their importance to the improvement and development of synthetic manufacturing processes and formulation stages. Deepening the concept of structure-reactivity and properties of materials."
$full.Font.Italic = 1

# --- Paragraph 14: Programa PT -> becomes single sentence (old Metodo value) ---
$p = $d.Paragraphs.Item(14)
$full = $d.Range($p.Range.Start, $p.Range.End - 1)
$full.Text = "Exposição e discussão de artigos, e desenvolvimento de experimentos propostos."

# --- Paragraph 17: Avaliacao ListBullet (Metodo/Criterio/Norma) ---
# Rebuild the run values in place via targeted Find & Replace, preserving the bold labels.
# Processed Norma -> Criterio -> Metodo (right-to-left) to avoid search-text collisions.
$p = $d.Paragraphs.Item(17)

$pRange = $p.Range
$f = $pRange.Find
$f.Execute("Aos alunos que tiverem freqüência mínima de 70% e média final menor que 5,0 e igual ou maior que 3,0, será dada recuperação com uma avaliação escrita. A média dessa avaliação somada com a média anterior das P1 e P2, se superior a cinco (5,0), levará a aprovação do aluno.", $true, $false, $false, $false, $false, $true, 1, $false, "Lehn, J. (1993). `"Supramolecular chemistry`". Science. 260 (5115): 1762–3. `v `vLehn, J.-M. (1995) Supramolecular Chemistry. Wiley-VCH. ISBN 978-3-527-29311-7`v`vNicholas J. Turro, V. Ramamurthy, J.C. Scaiano. Modern Molecular Photochemistry of Organic Molecules.`vSBN 978-1-891389-25-2, 1110 pages, Copyright 2010, Casebound.`v`vSilverstein, Robert M.; Webster Francis X.; Kiemle David J. Identificação Espectrométrica de Compostos Orgânicos, 7ª edição LTC", 2) | Out-Null

$pRange = $p.Range
$f = $pRange.Find
$f.Execute("Avaliação de seminários ministrados e da elaboração dos experimentos.", $true, $false, $false, $false, $false, $true, 1, $false, "Aos alunos que tiverem freqüência mínima de 70% e média final menor que 5,0 e igual ou maior que 3,0, será dada recuperação com uma avaliação escrita. A média dessa avaliação somada com a média anterior das P1 e P2, se superior a cinco (5,0), levará a aprovação do aluno.`v", 2) | Out-Null

$pRange = $p.Range
$f = $pRange.Find
$f.Execute("Exposição e discussão de artigos, e desenvolvimento de experimentos propostos.", $true, $false, $false, $false, $false, $true, 1, $false, "Avaliação de seminários ministrados e da elaboração dos experimentos.`v", 2) | Out-Null

# --- Paragraph 19: Bibliografia content -> becomes Docente bullet text ---
$p = $d.Paragraphs.Item(19)
$full = $d.Range($p.Range.Start, $p.Range.End - 1)
$full.Text = "210064 - Eduardo Rezende Triboni"

Write-Output "done"
